# Generate Report for Handback
# Update the "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" /
# "Correspond Handback DateTime" timestamps that were refreshed when the
# handback report was regenerated.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$overview.Range("G4").Value = "2016-09-05 20:53:13"

$zhcn.Range("H4").Value = "2016-09-05 20:53:06"
$zhcn.Range("K4").Value = "2016-09-05 20:53:35"

$dede.Range("H4").Value = "2016-09-05 20:53:13"
$dede.Range("K4").Value = "2016-09-05 20:53:43"
